$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything first so the shared-string table rebuilds from scratch
$ws.Cells.ClearContents()

# Write header row (row 1) - establishes shared strings 0-19
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# Write column A (rows 2-7) first so FAPs/MuSCs register before Ntrk3/Ptprs/ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "MuSCs"
$ws.Range("A6").Value = "MuSCs"
$ws.Range("A7").Value = "MuSCs"

# Column B (rows 2-7) -> registers Ntrk3
$ws.Range("B2").Value = "Ntrk3"
$ws.Range("B3").Value = "Ntrk3"
$ws.Range("B4").Value = "Ntrk3"
$ws.Range("B5").Value = "Ntrk3"
$ws.Range("B6").Value = "Ntrk3"
$ws.Range("B7").Value = "Ntrk3"

# Column C (rows 2-7) -> registers Ptprs
$ws.Range("C2").Value = "Ptprs"
$ws.Range("C3").Value = "Ptprs"
$ws.Range("C4").Value = "Ptprs"
$ws.Range("C5").Value = "Ptprs"
$ws.Range("C6").Value = "Ptprs"
$ws.Range("C7").Value = "Ptprs"

# Column D (rows 2-7) -> registers ECs last
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"

# Numeric columns E-T (rows 2-7)
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 2.468673666666667
$ws.Range("H2").Value = 7.406021
$ws.Range("I2").Value = 0.635345274347677
$ws.Range("J2").Value = 0.635345274347677
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 4.099037598280333
$ws.Range("R2").Value = 36.891338384523
$ws.Range("S2").Value = 0.0236017413541709
$ws.Range("T2").Value = 0.0236017413541709

$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 2.468673666666667
$ws.Range("H3").Value = 7.406021
$ws.Range("I3").Value = 0.635345274347677
$ws.Range("J3").Value = 0.635345274347677
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 25.17096033333333
$ws.Range("N3").Value = 75.512881
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("Q3").Value = 62.13888693961121
$ws.Range("R3").Value = 559.2499824565009
$ws.Range("S3").Value = 0.3577878715238055
$ws.Range("T3").Value = 0.3577878715238056

$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 2.468673666666667
$ws.Range("H4").Value = 7.406021
$ws.Range("I4").Value = 0.635345274347677
$ws.Range("J4").Value = 0.635345274347677
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 44.10580512003167
$ws.Range("R4").Value = 396.952246080285
$ws.Range("S4").Value = 0.2539556614697004
$ws.Range("T4").Value = 0.2539556614697004

$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 1.416888666666667
$ws.Range("H5").Value = 4.250666
$ws.Range("I5").Value = 0.364654725652323
$ws.Range("J5").Value = 0.364654725652323
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 1.660421
$ws.Range("N5").Value = 4.981262999999999
$ws.Range("O5").Value = 0.03714789785507311
$ws.Range("P5").Value = 0.03714789785507311
$ws.Range("Q5").Value = 2.352631696795333
$ws.Range("R5").Value = 21.173685271158
$ws.Range("S5").Value = 0.0135461565009022
$ws.Range("T5").Value = 0.0135461565009022

$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 1.416888666666667
$ws.Range("H6").Value = 4.250666
$ws.Range("I6").Value = 0.364654725652323
$ws.Range("J6").Value = 0.364654725652323
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 25.17096033333333
$ws.Range("N6").Value = 75.512881
$ws.Range("O6").Value = 0.5631392661118858
$ws.Range("P6").Value = 0.5631392661118859
$ws.Range("Q6").Value = 35.66444842541622
$ws.Range("R6").Value = 320.9800358287459
$ws.Range("S6").Value = 0.2053513945880802
$ws.Range("T6").Value = 0.2053513945880803

$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 1.416888666666667
$ws.Range("H7").Value = 4.250666
$ws.Range("I7").Value = 0.364654725652323
$ws.Range("J7").Value = 0.364654725652323
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 17.866195
$ws.Range("N7").Value = 53.598585
$ws.Range("O7").Value = 0.399712836033041
$ws.Range("P7").Value = 0.399712836033041
$ws.Range("Q7").Value = 25.31440921195667
$ws.Range("R7").Value = 227.82968290761
$ws.Range("S7").Value = 0.1457571745633405
$ws.Range("T7").Value = 0.1457571745633405

